$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sites")
$ws = $wb.Worksheets.Item("Articles")

# --- Row 19 / Row 20: fill in the SiteId so the VLOOKUP resolves to "pl" instead of #N/A ---
$ws.Range("A19").Value = 10
$ws.Range("A20").Value = 10

# --- Add the three new "Builds" hyperlinks (column C) for the pl rows/queue row ---
# Order matters: it controls the order new entries land in the shared string table,
# matching how the workbook was actually authored (Ornn, then Nautilus, then Sett).
$ws.Hyperlinks.Add($ws.Range("C19"), "https://buildpl.wordpress.com/2021/01/08/ornn-build-i-renekton-build-refleksja-na-top-lane/", "", "", "https://buildpl.wordpress.com/2021/01/08/ornn-build-i-renekton-build-refleksja-na-top-lane/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "https://build-pt.tumblr.com/post/639740051705708544/nautilus-build-yuumi-build-mad-kaiser-e-rge?is_related_post=1", "", "", "https://build-pt.tumblr.com/post/639740051705708544/nautilus-build-yuumi-build-mad-kaiser-e-rge?is_related_post=1") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C20"), "https://buildpl.wordpress.com/2021/01/10/sett-build-refleksja-na-junglerzy/", "", "", "https://buildpl.wordpress.com/2021/01/10/sett-build-refleksja-na-junglerzy/") | Out-Null

# Hyperlinks.Add() re-styles the cell with its own font variant; re-apply the
# same "Builds" link formatting already used by the other link cells (e.g. C16)
# so the new cells share the existing style instead of a near-duplicate one.
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)

# --- Added dates (column E): these rows moved from "In Queue" / "not posted yet" to real dates ---
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E17").Value = 44204
$ws.Range("E19").Value = 44204
$ws.Range("E20").Value = 44206

$excel.CutCopyMode = $false

# --- View/selection refresh to match where the author last clicked ---
$ws1.Select()
$ws1.Range("D11").Select() | Out-Null
$ws.Select()
$ws.Range("C28").Select() | Out-Null
